$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.2594103048008066
$ws.Range("J2").Value = 0.2594103048008066
$ws.Range("M2").Value = 2.113523666666667
$ws.Range("N2").Value = 6.340571000000001
$ws.Range("O2").Value = 0.2651220308693004
$ws.Range("P2").Value = 0.2651220308693004
$ws.Range("Q2").Value = 0.02827683313633334
$ws.Range("R2").Value = 0.254491498227
$ws.Range("S2").Value = 0.06877538683721407
$ws.Range("T2").Value = 0.06877538683721406

$ws.Range("I3").Value = 0.2594103048008066
$ws.Range("J3").Value = 0.2594103048008066
$ws.Range("O3").Value = 0.2869289465860668
$ws.Range("P3").Value = 0.2869289465860668
$ws.Range("S3").Value = 0.07443232549006594
$ws.Range("T3").Value = 0.07443232549006593

$ws.Range("I4").Value = 0.2594103048008066
$ws.Range("J4").Value = 0.2594103048008066
$ws.Range("M4").Value = 1.164746666666667
$ws.Range("N4").Value = 3.49424
$ws.Range("O4").Value = 0.1461067158059967
$ws.Range("P4").Value = 0.1461067158059966
$ws.Range("Q4").Value = 0.01558314565333333
$ws.Range("R4").Value = 0.14024831088
$ws.Range("S4").Value = 0.03790158768067842
$ws.Range("T4").Value = 0.03790158768067841

$ws.Range("I5").Value = 0.2594103048008066
$ws.Range("J5").Value = 0.2594103048008066
$ws.Range("M5").Value = 2.406253666666667
$ws.Range("N5").Value = 7.218761
$ws.Range("O5").Value = 0.3018423067386362
$ws.Range("P5").Value = 0.3018423067386362
$ws.Range("Q5").Value = 0.03219326780633334
$ws.Range("R5").Value = 0.289739410257
$ws.Range("S5").Value = 0.07830100479284818
$ws.Range("T5").Value = 0.07830100479284816

$ws.Range("G6").Value = 0.03819566666666666
$ws.Range("H6").Value = 0.114587
$ws.Range("I6").Value = 0.7405896951991934
$ws.Range("J6").Value = 0.7405896951991934
$ws.Range("M6").Value = 2.113523666666667
$ws.Range("N6").Value = 6.340571000000001
$ws.Range("O6").Value = 0.2651220308693004
$ws.Range("P6").Value = 0.2651220308693004
$ws.Range("Q6").Value = 0.08072744546411112
$ws.Range("R6").Value = 0.726547009177
$ws.Range("S6").Value = 0.1963466440320863
$ws.Range("T6").Value = 0.1963466440320863

$ws.Range("G7").Value = 0.03819566666666666
$ws.Range("H7").Value = 0.114587
$ws.Range("I7").Value = 0.7405896951991934
$ws.Range("J7").Value = 0.7405896951991934
$ws.Range("O7").Value = 0.2869289465860668
$ws.Range("P7").Value = 0.2869289465860668
$ws.Range("Q7").Value = 0.08736746928066666
$ws.Range("R7").Value = 0.7863072235259999
$ws.Range("S7").Value = 0.2124966210960008
$ws.Range("T7").Value = 0.2124966210960008

$ws.Range("G8").Value = 0.03819566666666666
$ws.Range("H8").Value = 0.114587
$ws.Range("I8").Value = 0.7405896951991934
$ws.Range("J8").Value = 0.7405896951991934
$ws.Range("M8").Value = 1.164746666666667
$ws.Range("N8").Value = 3.49424
$ws.Range("O8").Value = 0.1461067158059967
$ws.Range("P8").Value = 0.1461067158059966
$ws.Range("Q8").Value = 0.04448827543111111
$ws.Range("R8").Value = 0.40039447888
$ws.Range("S8").Value = 0.1082051281253182
$ws.Range("T8").Value = 0.1082051281253182

$ws.Range("G9").Value = 0.03819566666666666
$ws.Range("H9").Value = 0.114587
$ws.Range("I9").Value = 0.7405896951991934
$ws.Range("J9").Value = 0.7405896951991934
$ws.Range("M9").Value = 2.406253666666667
$ws.Range("N9").Value = 7.218761
$ws.Range("O9").Value = 0.3018423067386362
$ws.Range("P9").Value = 0.3018423067386362
$ws.Range("Q9").Value = 0.09190846296744444
$ws.Range("R9").Value = 0.8271761667069999
$ws.Range("S9").Value = 0.223541301945788
$ws.Range("T9").Value = 0.223541301945788
